$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83; all existing rows from 83 downward
# shift down by one (old row 83 -> 84, ..., old row 184 -> 185).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly record.
$ws.Range("A83").Value = 8
$ws.Range("B83").Value = "Terminal La Palmera de La Serena"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44895
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 100112040
$ws.Range("G83").Value = "Cilantro"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 2000
$ws.Range("M83").Value = 1750
$ws.Range("N83").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O83").Value = 'Provincia del Elquí'
$ws.Range("P83").Value = 1167
$ws.Range("Q83").Value = 1.5
$ws.Range("R83").Value = "Hortaliza"
